# Update weekly fruit/vegetable price data: rotate the values of rows 2, 3, 4
# and 9 (columns D, J, K, L, M, P) so the records are reordered chronologically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the four affected rows so we can
# reassign them without clobbering data we still need to read.
$rows = @(2, 3, 4, 9)
$cols = @("D", "J", "K", "L", "M", "P")

$original = @{}
foreach ($r in $rows) {
    $original[$r] = @{}
    foreach ($c in $cols) {
        $original[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row order: row2 <- old row9, row3 <- old row2, row4 <- old row3, row9 <- old row4
$mapping = @{
    2 = 9
    3 = 2
    4 = 3
    9 = 4
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $original[$srcRow][$c]
    }
}
